# profile onboarding scripts ENW to Neon
# Adds three new test-case rows (Profile61, Profile62, Profile63) to the
# "Test Cases" sheet, covering ENW -> Neon on-boarding welcome modal
# scenarios for the "Project Neon", "Profile" and "Account" links.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 62:64 from the formatting of the last existing data row (61)
# so the new rows pick up the same borders/fonts/column styles.
$ws.Range("A61:E61").Copy()
$ws.Range("A62:E64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 62 - Profile61 / OPQA-2087 - "Project Neon" link
$ws.Range("A62").Value = "Profile61"
$ws.Range("B62").Value = "OPQA-2087"
$ws.Range("C62").Value = 'Verify that ENW user who has not been on-boarded to Neon when clicks on "Project Neon" link from within ENW, shall be presented with the Neon on-boarding welcome modal.'
$ws.Range("D62").Value = "Y"

# Row 63 - Profile62 / OPQA-2089 - "Profile" link
$ws.Range("A63").Value = "Profile62"
$ws.Range("B63").Value = "OPQA-2089"
$ws.Range("C63").Value = 'Verify that ENW user who has not been on-boarded to Neon when clicks on "Profile" link from within ENW, shall be presented with the Neon on-boarding welcome modal.'
$ws.Range("D63").Value = "Y"

# Row 64 - Profile63 / OPQA-2090 - "Account" link
$ws.Range("A64").Value = "Profile63"
$ws.Range("B64").Value = "OPQA-2090"
$ws.Range("C64").Value = 'Verify that ENW user who has not been on-boarded to Neon when clicks on "Account" link from within ENW, shall be presented with the Neon on-boarding welcome modal.'
$ws.Range("D64").Value = "Y"

# Update the current selection to match the edited workbook state.
$ws.Range("C45").Select()
